# Auto-generated edit script: updates crypto price/volume table (rows 2-51)
# Applies the Tue Jan  2 01:09:42 UTC 2024 data refresh from GitHub Actions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, $Address, $Text)
    # Leading apostrophe forces Excel to treat the value as literal text,
    # even when it looks numeric (e.g. '45.166.33' or '0.633') or carries
    # significant whitespace/percent formatting (e.g. '  +6.39%  ').
    $Sheet.Range($Address).Value = "'" + $Text
    # Resetting the style back to Normal drops the incidental text-format
    # stamp so the cell keeps the workbook's default (unstyled) look.
    $Sheet.Range($Address).Style = "Normal"
}

# Row 2
Set-TextCell $ws 'D2' '45.166.33'
Set-TextCell $ws 'E2' '  +6.39%  '

# Row 3
Set-TextCell $ws 'D3' '2.390.98'
Set-TextCell $ws 'E3' '  +4.18%  '

# Row 4
Set-TextCell $ws 'E4' '  -0.70%  '

# Row 5
Set-TextCell $ws 'B5' 'Solana'
Set-TextCell $ws 'C5' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell $ws 'D5' '111.69'
Set-TextCell $ws 'E5' '  +8.90%  '

# Row 6
Set-TextCell $ws 'B6' 'BNB'
Set-TextCell $ws 'C6' 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextCell $ws 'D6' '318.43'
Set-TextCell $ws 'E6' '  +1.15%  '

# Row 7
Set-TextCell $ws 'D7' '0.633'
Set-TextCell $ws 'E7' '  +1.55%  '

# Row 8
Set-TextCell $ws 'E8' '  -0.38%  '

# Row 9
Set-TextCell $ws 'D9' '0.636'
Set-TextCell $ws 'E9' '  +5.75%  '

# Row 10
Set-TextCell $ws 'D10' '42.65'
Set-TextCell $ws 'E10' '  +9.32%  '

# Row 11
Set-TextCell $ws 'D11' '0.0932'
Set-TextCell $ws 'E11' '  +3.51%  '

# Row 12
Set-TextCell $ws 'D12' '8.72'
Set-TextCell $ws 'E12' '  +5.09%  '

# Row 13
Set-TextCell $ws 'E13' '  +3.69%  '

# Row 14
Set-TextCell $ws 'E14' '  -0.53%  '

# Row 15
Set-TextCell $ws 'D15' '15.79'
Set-TextCell $ws 'E15' '  +4.26%  '

# Row 16
Set-TextCell $ws 'D16' '2.745.01'
Set-TextCell $ws 'E16' '  +3.77%  '

# Row 17
Set-TextCell $ws 'D17' '2.393.37'
Set-TextCell $ws 'E17' '  +4.15%  '

# Row 18
Set-TextCell $ws 'D18' '45.091.94'
Set-TextCell $ws 'E18' '  +5.78%  '

# Row 19
Set-TextCell $ws 'D19' '7.66'
Set-TextCell $ws 'E19' '  +4.90%  '

# Row 20
Set-TextCell $ws 'D20' '0.0000108'
Set-TextCell $ws 'E20' '  +3.37%  '

# Row 21
Set-TextCell $ws 'D21' '13.15'
Set-TextCell $ws 'E21' '  -3.64%  '

# Row 22
Set-TextCell $ws 'D22' '75.57'
Set-TextCell $ws 'E22' '  +3.47%  '

# Row 23
Set-TextCell $ws 'D23' '3.58'
Set-TextCell $ws 'E23' '  +2.59%  '

# Row 24
Set-TextCell $ws 'D24' '269.93'
Set-TextCell $ws 'E24' '  +2.76%  '

# Row 25
Set-TextCell $ws 'D25' '2.32'
Set-TextCell $ws 'E25' '  +6.34%  '

# Row 26
Set-TextCell $ws 'E26' '  -0.58%  '

# Row 27
Set-TextCell $ws 'D27' '7.70'
Set-TextCell $ws 'E27' '  +13.39%  '

# Row 28
Set-TextCell $ws 'D28' '11.32'
Set-TextCell $ws 'E28' '  +5.63%  '

# Row 29
Set-TextCell $ws 'D29' '2.33'
Set-TextCell $ws 'E29' '  -0.38%  '

# Row 30
Set-TextCell $ws 'D30' '39.45'
Set-TextCell $ws 'E30' '  +9.82%  '

# Row 31
Set-TextCell $ws 'D31' '22.87'
Set-TextCell $ws 'E31' '  +2.54%  '

# Row 32
Set-TextCell $ws 'D32' '169.65'
Set-TextCell $ws 'E32' '  +2.18%  '

# Row 33
Set-TextCell $ws 'D33' '0.0921'
Set-TextCell $ws 'E33' '  +6.37%  '

# Row 34
Set-TextCell $ws 'D34' '3.00'
Set-TextCell $ws 'E34' '  +17.15%  '

# Row 35
Set-TextCell $ws 'D35' '0.133'
Set-TextCell $ws 'E35' '  +2.21%  '

# Row 36
Set-TextCell $ws 'D36' '0.120'
Set-TextCell $ws 'E36' '  +5.83%  '

# Row 37
Set-TextCell $ws 'D37' '4.82'
Set-TextCell $ws 'E37' '  +6.20%  '

# Row 38
Set-TextCell $ws 'D38' '0.0367'
Set-TextCell $ws 'E38' '  +5.35%  '

# Row 39
Set-TextCell $ws 'D39' '2.96'
Set-TextCell $ws 'E39' '  +10.81%  '

# Row 40
Set-TextCell $ws 'D40' '3.90'
Set-TextCell $ws 'E40' '  +3.66%  '

# Row 41
Set-TextCell $ws 'D41' '1.76'
Set-TextCell $ws 'E41' '  +11.40%  '

# Row 42
Set-TextCell $ws 'D42' '104.95'
Set-TextCell $ws 'E42' '  +7.82%  '

# Row 43
Set-TextCell $ws 'D43' '13.89'
Set-TextCell $ws 'E43' '  +14.94%  '

# Row 44
Set-TextCell $ws 'D44' '0.240'
Set-TextCell $ws 'E44' '  +5.92%  '

# Row 45
Set-TextCell $ws 'D45' '72.03'
Set-TextCell $ws 'E45' '  +3.53%  '

# Row 46
Set-TextCell $ws 'D46' '0.999'
Set-TextCell $ws 'E46' '  -0.56%  '

# Row 47
Set-TextCell $ws 'D47' '118.35'
Set-TextCell $ws 'E47' '  +6.69%  '

# Row 48
Set-TextCell $ws 'D48' '80.46'
Set-TextCell $ws 'E48' '  +0.10%  '

# Row 49
Set-TextCell $ws 'B49' 'THORChain'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextCell $ws 'D49' '5.50'
Set-TextCell $ws 'E49' '  +6.17%  '

# Row 50
Set-TextCell $ws 'B50' 'FraxShare'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws 'D50' '9.08'
Set-TextCell $ws 'E50' '  +4.77%  '

# Row 51
Set-TextCell $ws 'B51' 'TheGraph'
Set-TextCell $ws 'C51' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell $ws 'D51' '0.218'
Set-TextCell $ws 'E51' '  +14.76%  '
